$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Pt7"

$ws.Range("A1").Value = "Ieeg name"
$ws.Range("B1").Value = "Current"
$ws.Range("C1").Value = "Main stim start time"
$ws.Range("D1").Value = "Main stim end time"
$ws.Range("E1").Value = "Time breaks"
$ws.Range("F1").Value = "Electrodes"
$ws.Range("G1").Value = "Clinical Stim seizure elecs"
$ws.Range("H1").Value = "Suspected SOZ anatomic"
$ws.Range("I1").Value = "Afterdischarges"
$ws.Range("J1").Value = "Clinical effects"
$ws.Range("K1").Value = "Seizures"
$ws.Range("L1").Value = "Current test electrodes"
$ws.Range("M1").Value = "Electrode"
$ws.Range("N1").Value = "Anatomical target"
$ws.Range("O1").Value = "Other"
$ws.Range("P1").Value = "visuallyBadChannels"
$ws.Range("A2").Value = "HUP218_CCEP"
$ws.Range("B2").Value = "3 mA"
$ws.Range("D2").Value = "end"
$ws.Range("F2").Value = "RA1"
$ws.Range("G2").Value = "RA1"
$ws.Range("H2").Value = "mesial temporal"
$ws.Range("I2").Value = "none"
$ws.Range("J2").Value = "none"
$ws.Range("K2").Value = "none"
$ws.Range("M2").Value = "RA"
$ws.Range("N2").Value = "right amygdala"
$ws.Range("F3").Value = "RA2"
$ws.Range("G3").Value = "RB1"
$ws.Range("M3").Value = "RB"
$ws.Range("N3").Value = "right anterior hippocampus"
$ws.Range("F4").Value = "RA3"
$ws.Range("M4").Value = "RC"
$ws.Range("N4").Value = "right posterior hippocampus"
$ws.Range("F5").Value = "RA6"
$ws.Range("M5").Value = "RD"
$ws.Range("N5").Value = "right anterior insula"
$ws.Range("F6").Value = "RA7"
$ws.Range("M6").Value = "RE"
$ws.Range("N6").Value = "right posterior insula"
$ws.Range("F7").Value = "RA8"
$ws.Range("M7").Value = "RF"
$ws.Range("N7").Value = "right anterior cingulate"
$ws.Range("F8").Value = "RA9"
$ws.Range("M8").Value = "RG"
$ws.Range("N8").Value = "right mid cingulate"
$ws.Range("F9").Value = "RB1"
$ws.Range("M9").Value = "RH"
$ws.Range("N9").Value = "right orbitofrontal/PET"
$ws.Range("F10").Value = "RB2"
$ws.Range("M10").Value = "RI"
$ws.Range("N10").Value = "right superior frontal gyrus"
$ws.Range("F11").Value = "RB3"
$ws.Range("M11").Value = "RJ"
$ws.Range("N11").Value = "right SSMA"
$ws.Range("F12").Value = "RC1"
$ws.Range("M12").Value = "RK"
$ws.Range("N12").Value = "right TPO"
$ws.Range("F13").Value = "RC2"
$ws.Range("M13").Value = "RL"
$ws.Range("N13").Value = "right basal temporal"
$ws.Range("F14").Value = "RC8"
$ws.Range("M14").Value = "RM"
$ws.Range("N14").Value = "right anterior insula revision"
$ws.Range("F15").Value = "RD1"
$ws.Range("M15").Value = "RN"
$ws.Range("N15").Value = "right anterior temporal revision"
$ws.Range("F16").Value = "RD2"
$ws.Range("M16").Value = "RO"
$ws.Range("N16").Value = "right orbitofrontal revision"
$ws.Range("F17").Value = "RD3"
$ws.Range("M17").Value = "RP"
$ws.Range("N17").Value = "right mesial frontal revision"
$ws.Range("F18").Value = "RD4"
$ws.Range("M18").Value = "RR"
$ws.Range("N18").Value = "right frontal pole revision"
$ws.Range("F19").Value = "RD5"
$ws.Range("F20").Value = "RD6"
$ws.Range("F21").Value = "RE1"
$ws.Range("F22").Value = "RE2"
$ws.Range("F23").Value = "RE3"
$ws.Range("F24").Value = "RE4"
$ws.Range("F25").Value = "RE5"
$ws.Range("F26").Value = "RG1"
$ws.Range("F27").Value = "RG2"
$ws.Range("F28").Value = "RG4"
$ws.Range("F29").Value = "RG5"
$ws.Range("F30").Value = "RG6"
$ws.Range("F31").Value = "RG7"
$ws.Range("F32").Value = "RI1"
$ws.Range("F33").Value = "RI2"
$ws.Range("F34").Value = "RI3"
$ws.Range("F35").Value = "RI4"
$ws.Range("F36").Value = "RI5"
$ws.Range("F37").Value = "RK2"
$ws.Range("F38").Value = "RK3"
$ws.Range("F39").Value = "RK4"
$ws.Range("F40").Value = "RL1"
$ws.Range("F41").Value = "RL3"
$ws.Range("F42").Value = "RL7"
$ws.Range("F43").Value = "RL9"
$ws.Range("F44").Value = "RM1"
$ws.Range("F45").Value = "RM2"
$ws.Range("F46").Value = "RM3"
$ws.Range("F47").Value = "RM4"
$ws.Range("F48").Value = "RM5"
$ws.Range("F49").Value = "RM6"
$ws.Range("F50").Value = "RM7"
$ws.Range("F51").Value = "RN1"
$ws.Range("F52").Value = "RN2"
$ws.Range("F53").Value = "RN3"
$ws.Range("F54").Value = "RN4"
$ws.Range("F55").Value = "RN7"
$ws.Range("F56").Value = "RP1"
$ws.Range("F57").Value = "RP2"
$ws.Range("F58").Value = "RP3"
$ws.Range("F59").Value = "RR1"
$ws.Range("F60").Value = "RR2"
$ws.Range("F61").Value = "RR4"
$ws.Range("F62").Value = "RR5"
$ws.Range("F63").Value = "RR6"
$ws.Range("F64").Value = "RR7"

$ws.Range("C7").Select()
$wb.Worksheets.Item(7).Activate()
